$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("<one>",      "<one>",      9),
    @("<oil>",      "<oil>",      0),
    @("<up>",       "<up>",       10),
    @("<up>",       "<up>",       13),
    @("<like>",     "<like>",     10),
    @("<and>",      "<and>",      10),
    @("<shift>",    "<shift>",    6),
    @("<said>",     "<said>",     11),
    @("<tango>",    "<tango>",    11),
    @("<mike>",     "<make>",     8),
    @("<kilo>",     "<kilo>",     6),
    @("<a>",        "<a>",        13),
    @("<may>",      "<may>",      11),
    @("<water>",    "<water>",    11),
    @("<zulu>",     "<zulu>",     6),
    @("<water>",    "<water>",    15),
    @("<uniform>",  "<uniform>",  8)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
